# Auto-generated script to update Leve profit calculation sheets
# with refreshed market-price data (scheduled runner sync).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3108.0356  # H64: 3211.7307 -> 3108.0356
$ws.Cells.Item(64, 10).Value = 3232.8096  # J64: 3387.842 -> 3232.8096
$ws.Cells.Item(64, 12).Value = 3232.8096  # L64: 3387.842 -> 3232.8096
$ws.Cells.Item(64, 14).Value = -3728.8096  # N64: -3883.842 -> -3728.8096
$ws.Cells.Item(67, 8).Value = 3108.0356  # H67: 3211.7307 -> 3108.0356
$ws.Cells.Item(67, 10).Value = 3232.8096  # J67: 3387.842 -> 3232.8096
$ws.Cells.Item(67, 12).Value = 3232.8096  # L67: 3387.842 -> 3232.8096
$ws.Cells.Item(67, 14).Value = -4948.809600000001  # N67: -5103.842000000001 -> -4948.809600000001
$ws.Cells.Item(74, 8).Value = 5271.8184  # H74: 4387.875 -> 5271.8184
$ws.Cells.Item(74, 9).Value = 3858.5715  # I74: 3352.1667 -> 3858.5715
$ws.Cells.Item(74, 10).Value = 7745  # J74: 7495 -> 7745
$ws.Cells.Item(74, 11).Value = 3858.5715  # K74: 3352.1667 -> 3858.5715
$ws.Cells.Item(74, 12).Value = 7745  # L74: 7495 -> 7745
$ws.Cells.Item(74, 13).Value = -2922.5715  # M74: -2416.1667 -> -2922.5715
$ws.Cells.Item(74, 14).Value = -9617  # N74: -9367 -> -9617
$ws.Cells.Item(77, 8).Value = 5271.8184  # H77: 4387.875 -> 5271.8184
$ws.Cells.Item(77, 9).Value = 3858.5715  # I77: 3352.1667 -> 3858.5715
$ws.Cells.Item(77, 10).Value = 7745  # J77: 7495 -> 7745
$ws.Cells.Item(77, 11).Value = 19292.8575  # K77: 16760.8335 -> 19292.8575
$ws.Cells.Item(77, 12).Value = 38725  # L77: 37475 -> 38725
$ws.Cells.Item(77, 13).Value = -14612.8575  # M77: -12080.8335 -> -14612.8575
$ws.Cells.Item(77, 14).Value = -48085  # N77: -46835 -> -48085
$ws.Cells.Item(100, 8).Value = 2224.1667  # H100: 2491.4285 -> 2224.1667
$ws.Cells.Item(100, 9).Value = 1667.0834  # I100: 1933.3334 -> 1667.0834
$ws.Cells.Item(100, 10).Value = 3338.3333  # J100: 3496 -> 3338.3333
$ws.Cells.Item(100, 11).Value = 1667.0834  # K100: 1933.3334 -> 1667.0834
$ws.Cells.Item(100, 12).Value = 3338.3333  # L100: 3496 -> 3338.3333
$ws.Cells.Item(100, 13).Value = -1126.0834  # M100: -1392.3334 -> -1126.0834
$ws.Cells.Item(100, 14).Value = -4420.3333  # N100: -4578 -> -4420.3333
$ws.Cells.Item(113, 8).Value = 3106.72  # H113: 3212.652 -> 3106.72
$ws.Cells.Item(113, 9).Value = 2600.5  # I113: 2591.2727 -> 2600.5
$ws.Cells.Item(113, 10).Value = 3444.2  # J113: 3782.25 -> 3444.2
$ws.Cells.Item(113, 11).Value = 2600.5  # K113: 2591.2727 -> 2600.5
$ws.Cells.Item(113, 12).Value = 3444.2  # L113: 3782.25 -> 3444.2
$ws.Cells.Item(113, 13).Value = 653.5  # M113: 662.7273 -> 653.5
$ws.Cells.Item(113, 14).Value = -9952.200000000001  # N113: -10290.25 -> -9952.200000000001
$ws.Cells.Item(129, 8).Value = 1369.878  # H129: 1298.1957 -> 1369.878
$ws.Cells.Item(129, 10).Value = 1676.8  # J129: 1538.7428 -> 1676.8
$ws.Cells.Item(129, 12).Value = 5030.4  # L129: 4616.2284 -> 5030.4
$ws.Cells.Item(129, 14).Value = -15030.4  # N129: -14616.2284 -> -15030.4
$ws.Cells.Item(132, 8).Value = 6609.174  # H132: 6560.3125 -> 6609.174
$ws.Cells.Item(132, 9).Value = 5582.9697  # I132: 5568.6 -> 5582.9697
$ws.Cells.Item(132, 10).Value = 9214.154  # J132: 9230.308000000001 -> 9214.154
$ws.Cells.Item(132, 11).Value = 16748.9091  # K132: 16705.8 -> 16748.9091
$ws.Cells.Item(132, 12).Value = 27642.462  # L132: 27690.924 -> 27642.462
$ws.Cells.Item(132, 13).Value = -14218.9091  # M132: -14175.8 -> -14218.9091
$ws.Cells.Item(132, 14).Value = -32702.462  # N132: -32750.924 -> -32702.462
$ws.Cells.Item(138, 8).Value = 1908.0312  # H138: 1816.425 -> 1908.0312
$ws.Cells.Item(138, 10).Value = 2147.5  # J138: 1837.5 -> 2147.5
$ws.Cells.Item(138, 12).Value = 6442.5  # L138: 5512.5 -> 6442.5
$ws.Cells.Item(138, 14).Value = -16722.5  # N138: -15792.5 -> -16722.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1635.3  # H45: 1692.84 -> 1635.3
$ws.Cells.Item(45, 9).Value = 1300.5625  # I45: 1423.75 -> 1300.5625
$ws.Cells.Item(45, 10).Value = 2017.8572  # J45: 1941.2307 -> 2017.8572
$ws.Cells.Item(45, 11).Value = 1300.5625  # K45: 1423.75 -> 1300.5625
$ws.Cells.Item(45, 12).Value = 2017.8572  # L45: 1941.2307 -> 2017.8572
$ws.Cells.Item(45, 13).Value = -923.5625  # M45: -1046.75 -> -923.5625
$ws.Cells.Item(45, 14).Value = -2771.8572  # N45: -2695.2307 -> -2771.8572
$ws.Cells.Item(97, 8).Value = 13548.6875  # H97: 23484.334 -> 13548.6875
$ws.Cells.Item(97, 9).Value = 16491.46  # I97: 26311.125 -> 16491.46
$ws.Cells.Item(97, 10).Value = 796.6667  # J97: 870 -> 796.6667
$ws.Cells.Item(97, 11).Value = 16491.46  # K97: 26311.125 -> 16491.46
$ws.Cells.Item(97, 12).Value = 796.6667  # L97: 870 -> 796.6667
$ws.Cells.Item(97, 13).Value = -15995.46  # M97: -25815.125 -> -15995.46
$ws.Cells.Item(97, 14).Value = -1788.6667  # N97: -1862 -> -1788.6667
$ws.Cells.Item(102, 8).Value = 3874.6667  # H102: 4000 -> 3874.6667
$ws.Cells.Item(102, 9).Value = 3517.7778  # I102: 0 -> 3517.7778
$ws.Cells.Item(102, 10).Value = 4410  # J102: 4000 -> 4410
$ws.Cells.Item(102, 11).Value = 3517.7778  # K102: 0 -> 3517.7778
$ws.Cells.Item(102, 12).Value = 4410  # L102: 4000 -> 4410
$ws.Cells.Item(102, 13).Value = -1895.7778  # M102: None -> -1895.7778
$ws.Cells.Item(102, 14).Value = -7654  # N102: -7244 -> -7654

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 984.7406999999999  # H86: 1142.2667 -> 984.7406999999999
$ws.Cells.Item(86, 9).Value = 968.5625  # I86: 1112 -> 968.5625
$ws.Cells.Item(86, 10).Value = 1008.2727  # J86: 1202.8 -> 1008.2727
$ws.Cells.Item(86, 11).Value = 968.5625  # K86: 1112 -> 968.5625
$ws.Cells.Item(86, 12).Value = 1008.2727  # L86: 1202.8 -> 1008.2727
$ws.Cells.Item(86, 13).Value = 154.4375  # M86: 11 -> 154.4375
$ws.Cells.Item(86, 14).Value = -3254.2727  # N86: -3448.8 -> -3254.2727
$ws.Cells.Item(89, 8).Value = 984.7406999999999  # H89: 1142.2667 -> 984.7406999999999
$ws.Cells.Item(89, 9).Value = 968.5625  # I89: 1112 -> 968.5625
$ws.Cells.Item(89, 10).Value = 1008.2727  # J89: 1202.8 -> 1008.2727
$ws.Cells.Item(89, 11).Value = 4842.8125  # K89: 5560 -> 4842.8125
$ws.Cells.Item(89, 12).Value = 5041.363499999999  # L89: 6014 -> 5041.363499999999
$ws.Cells.Item(89, 13).Value = 773.1875  # M89: 56 -> 773.1875
$ws.Cells.Item(89, 14).Value = -16273.3635  # N89: -17246 -> -16273.3635
$ws.Cells.Item(94, 8).Value = 10361.571  # H94: 1324.3334 -> 10361.571
$ws.Cells.Item(94, 9).Value = 11888.556  # I94: 1418.1666 -> 11888.556
$ws.Cells.Item(94, 10).Value = 1199.6666  # J94: 1136.6666 -> 1199.6666
$ws.Cells.Item(94, 11).Value = 11888.556  # K94: 1418.1666 -> 11888.556
$ws.Cells.Item(94, 12).Value = 1199.6666  # L94: 1136.6666 -> 1199.6666
$ws.Cells.Item(94, 13).Value = -11437.556  # M94: -967.1666 -> -11437.556
$ws.Cells.Item(94, 14).Value = -2101.6666  # N94: -2038.6666 -> -2101.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4890.25  # H16: 3137.3845 -> 4890.25
$ws.Cells.Item(16, 9).Value = 5187  # I16: 3208.6 -> 5187
$ws.Cells.Item(16, 10).Value = 4000  # J16: 2900 -> 4000
$ws.Cells.Item(16, 11).Value = 5187  # K16: 3208.6 -> 5187
$ws.Cells.Item(16, 12).Value = 4000  # L16: 2900 -> 4000
$ws.Cells.Item(16, 13).Value = -4900  # M16: -2921.6 -> -4900
$ws.Cells.Item(16, 14).Value = -4574  # N16: -3474 -> -4574
$ws.Cells.Item(25, 8).Value = 55259.75  # H25: 93346.336 -> 55259.75
$ws.Cells.Item(25, 10).Value = 55259.75  # J25: 93346.336 -> 55259.75
$ws.Cells.Item(25, 12).Value = 55259.75  # L25: 93346.336 -> 55259.75
$ws.Cells.Item(25, 14).Value = -55607.75  # N25: -93694.336 -> -55607.75
$ws.Cells.Item(62, 8).Value = 10058.077  # H62: 12609.8 -> 10058.077
$ws.Cells.Item(62, 9).Value = 2195  # I62: 2300 -> 2195
$ws.Cells.Item(62, 10).Value = 27750  # J62: 36666 -> 27750
$ws.Cells.Item(62, 11).Value = 2195  # K62: 2300 -> 2195
$ws.Cells.Item(62, 12).Value = 27750  # L62: 36666 -> 27750
$ws.Cells.Item(62, 13).Value = -1571  # M62: -1676 -> -1571
$ws.Cells.Item(62, 14).Value = -28998  # N62: -37914 -> -28998
$ws.Cells.Item(65, 8).Value = 10058.077  # H65: 12609.8 -> 10058.077
$ws.Cells.Item(65, 9).Value = 2195  # I65: 2300 -> 2195
$ws.Cells.Item(65, 10).Value = 27750  # J65: 36666 -> 27750
$ws.Cells.Item(65, 11).Value = 10975  # K65: 11500 -> 10975
$ws.Cells.Item(65, 12).Value = 138750  # L65: 183330 -> 138750
$ws.Cells.Item(65, 13).Value = -7855  # M65: -8380 -> -7855
$ws.Cells.Item(65, 14).Value = -144990  # N65: -189570 -> -144990
$ws.Cells.Item(105, 8).Value = 1060.25  # H105: 887 -> 1060.25
$ws.Cells.Item(105, 9).Value = 1060.25  # I105: 902.94446 -> 1060.25
$ws.Cells.Item(105, 10).Value = 0  # J105: 600 -> 0
$ws.Cells.Item(105, 11).Value = 1060.25  # K105: 902.94446 -> 1060.25
$ws.Cells.Item(105, 12).Value = 0  # L105: 600 -> 0
$ws.Cells.Item(105, 13).Value = 686.75  # M105: 844.05554 -> 686.75
$ws.Cells.Item(105, 14).ClearContents()  # N105: delete (was -4094)
$ws.Cells.Item(113, 8).Value = 4890.25  # H113: 3137.3845 -> 4890.25
$ws.Cells.Item(113, 9).Value = 5187  # I113: 3208.6 -> 5187
$ws.Cells.Item(113, 10).Value = 4000  # J113: 2900 -> 4000
$ws.Cells.Item(113, 11).Value = 5187  # K113: 3208.6 -> 5187
$ws.Cells.Item(113, 12).Value = 4000  # L113: 2900 -> 4000
$ws.Cells.Item(113, 13).Value = -3017  # M113: -1038.6 -> -3017
$ws.Cells.Item(113, 14).Value = -8340  # N113: -7240 -> -8340

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 545.81134  # H113: 655.19354 -> 545.81134
$ws.Cells.Item(113, 9).Value = 415.16666  # I113: 436.79166 -> 415.16666
$ws.Cells.Item(113, 10).Value = 1044.6364  # J113: 1404 -> 1044.6364
$ws.Cells.Item(113, 11).Value = 1245.49998  # K113: 1310.37498 -> 1245.49998
$ws.Cells.Item(113, 12).Value = 3133.9092  # L113: 4212 -> 3133.9092
$ws.Cells.Item(113, 13).Value = 924.5000199999999  # M113: 859.6250199999999 -> 924.5000199999999
$ws.Cells.Item(113, 14).Value = -7473.9092  # N113: -8552 -> -7473.9092
$ws.Cells.Item(122, 8).Value = 3503.8667  # H122: 3113.7354 -> 3503.8667
$ws.Cells.Item(122, 9).Value = 1800  # I122: 658.5 -> 1800
$ws.Cells.Item(122, 10).Value = 3562.6206  # J122: 3639.8572 -> 3562.6206
$ws.Cells.Item(122, 11).Value = 16200  # K122: 5926.5 -> 16200
$ws.Cells.Item(122, 12).Value = 32063.5854  # L122: 32758.7148 -> 32063.5854
$ws.Cells.Item(122, 13).Value = -13750  # M122: -3476.5 -> -13750
$ws.Cells.Item(122, 14).Value = -36963.5854  # N122: -37658.7148 -> -36963.5854

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2775.7693  # H80: 2648.5833 -> 2775.7693
$ws.Cells.Item(80, 9).Value = 2535.8  # I80: 2580 -> 2535.8
$ws.Cells.Item(80, 10).Value = 2925.75  # J80: 2717.1667 -> 2925.75
$ws.Cells.Item(80, 11).Value = 2535.8  # K80: 2580 -> 2535.8
$ws.Cells.Item(80, 12).Value = 2925.75  # L80: 2717.1667 -> 2925.75
$ws.Cells.Item(80, 13).Value = -1537.8  # M80: -1582 -> -1537.8
$ws.Cells.Item(80, 14).Value = -4921.75  # N80: -4713.1667 -> -4921.75
$ws.Cells.Item(83, 8).Value = 2775.7693  # H83: 2648.5833 -> 2775.7693
$ws.Cells.Item(83, 9).Value = 2535.8  # I83: 2580 -> 2535.8
$ws.Cells.Item(83, 10).Value = 2925.75  # J83: 2717.1667 -> 2925.75
$ws.Cells.Item(83, 11).Value = 12679  # K83: 12900 -> 12679
$ws.Cells.Item(83, 12).Value = 14628.75  # L83: 13585.8335 -> 14628.75
$ws.Cells.Item(83, 13).Value = -7687  # M83: -7908 -> -7687
$ws.Cells.Item(83, 14).Value = -24612.75  # N83: -23569.8335 -> -24612.75
$ws.Cells.Item(107, 8).Value = 849.0625  # H107: 772.5909 -> 849.0625
$ws.Cells.Item(107, 9).Value = 862.75  # I107: 717.6 -> 862.75
$ws.Cells.Item(107, 10).Value = 835.375  # J107: 890.4286 -> 835.375
$ws.Cells.Item(107, 11).Value = 862.75  # K107: 717.6 -> 862.75
$ws.Cells.Item(107, 12).Value = 835.375  # L107: 890.4286 -> 835.375
$ws.Cells.Item(107, 13).Value = 1057.25  # M107: 1202.4 -> 1057.25
$ws.Cells.Item(107, 14).Value = -4675.375  # N107: -4730.4286 -> -4675.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2268.182  # H68: 2412.5 -> 2268.182
$ws.Cells.Item(71, 8).Value = 2268.182  # H71: 2412.5 -> 2268.182
$ws.Cells.Item(93, 8).Value = 1756.4166  # H93: 2026.7222 -> 1756.4166
$ws.Cells.Item(93, 9).Value = 1357.2142  # I93: 1537.8 -> 1357.2142
$ws.Cells.Item(93, 10).Value = 2315.3  # J93: 2637.875 -> 2315.3
$ws.Cells.Item(93, 11).Value = 1357.2142  # K93: 1537.8 -> 1357.2142
$ws.Cells.Item(93, 12).Value = 2315.3  # L93: 2637.875 -> 2315.3
$ws.Cells.Item(93, 13).Value = -109.2141999999999  # M93: -289.8 -> -109.2141999999999
$ws.Cells.Item(93, 14).Value = -4811.3  # N93: -5133.875 -> -4811.3
$ws.Cells.Item(100, 8).Value = 2431.077  # H100: 3075.5 -> 2431.077
$ws.Cells.Item(100, 9).Value = 1116.6666  # I100: 1350 -> 1116.6666
$ws.Cells.Item(100, 10).Value = 3557.7144  # J100: 3650.6667 -> 3557.7144
$ws.Cells.Item(100, 11).Value = 1116.6666  # K100: 1350 -> 1116.6666
$ws.Cells.Item(100, 12).Value = 3557.7144  # L100: 3650.6667 -> 3557.7144
$ws.Cells.Item(100, 13).Value = -575.6666  # M100: -809 -> -575.6666
$ws.Cells.Item(100, 14).Value = -4639.7144  # N100: -4732.6667 -> -4639.7144
$ws.Cells.Item(132, 8).Value = 3018.6875  # H132: 32453.295 -> 3018.6875
$ws.Cells.Item(132, 9).Value = 2700.0557  # I132: 52775.75 -> 2700.0557
$ws.Cells.Item(132, 10).Value = 3428.3572  # J132: 3421.2144 -> 3428.3572
$ws.Cells.Item(132, 11).Value = 8100.1671  # K132: 158327.25 -> 8100.1671
$ws.Cells.Item(132, 12).Value = 10285.0716  # L132: 10263.6432 -> 10285.0716
$ws.Cells.Item(132, 13).Value = -5570.1671  # M132: -155797.25 -> -5570.1671
$ws.Cells.Item(132, 14).Value = -15345.0716  # N132: -15323.6432 -> -15345.0716
